# Update leveling-profit metrics on each job sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# following the latest market-board price pull from the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 16
$ws.Cells.Item(16, 8).Value = 12463.667  # H16
$ws.Cells.Item(16, 9).Value = 7549.25  # I16
$ws.Cells.Item(16, 11).Value = 7549.25  # K16
$ws.Cells.Item(16, 13).Value = -7319.25  # M16
# Row 21
$ws.Cells.Item(21, 8).Value = 5000  # H21
$ws.Cells.Item(21, 9).Value = 5000  # I21
$ws.Cells.Item(21, 11).Value = 5000  # K21
$ws.Cells.Item(21, 13).Value = -4532  # M21
# Row 23
$ws.Cells.Item(23, 8).Value = 5000  # H23
$ws.Cells.Item(23, 9).Value = 5000  # I23
$ws.Cells.Item(23, 11).Value = 5000  # K23
$ws.Cells.Item(23, 13).Value = -4766  # M23
# Row 64
$ws.Cells.Item(64, 8).Value = 6959.1304  # H64
$ws.Cells.Item(64, 9).Value = 5533.7  # I64
$ws.Cells.Item(64, 10).Value = 8055.615  # J64
$ws.Cells.Item(64, 11).Value = 5533.7  # K64
$ws.Cells.Item(64, 12).Value = 8055.615  # L64
$ws.Cells.Item(64, 13).Value = -5285.7  # M64
$ws.Cells.Item(64, 14).Value = -8551.615  # N64
# Row 67
$ws.Cells.Item(67, 8).Value = 6959.1304  # H67
$ws.Cells.Item(67, 9).Value = 5533.7  # I67
$ws.Cells.Item(67, 10).Value = 8055.615  # J67
$ws.Cells.Item(67, 11).Value = 5533.7  # K67
$ws.Cells.Item(67, 12).Value = 8055.615  # L67
$ws.Cells.Item(67, 13).Value = -4675.7  # M67
$ws.Cells.Item(67, 14).Value = -9771.615  # N67
# Row 76
$ws.Cells.Item(76, 8).Value = 18236.875  # H76
$ws.Cells.Item(76, 9).Value = 26479.2  # I76
$ws.Cells.Item(76, 10).Value = 4499.6665  # J76
$ws.Cells.Item(76, 11).Value = 26479.2  # K76
$ws.Cells.Item(76, 12).Value = 4499.6665  # L76
$ws.Cells.Item(76, 13).Value = -26164.2  # M76
$ws.Cells.Item(76, 14).Value = -5129.6665  # N76
# Row 79
$ws.Cells.Item(79, 8).Value = 18236.875  # H79
$ws.Cells.Item(79, 9).Value = 26479.2  # I79
$ws.Cells.Item(79, 10).Value = 4499.6665  # J79
$ws.Cells.Item(79, 11).Value = 26479.2  # K79
$ws.Cells.Item(79, 12).Value = 4499.6665  # L79
$ws.Cells.Item(79, 13).Value = -25387.2  # M79
$ws.Cells.Item(79, 14).Value = -6683.6665  # N79
# Row 137
$ws.Cells.Item(137, 8).Value = 6689.5415  # H137
$ws.Cells.Item(137, 9).Value = 8817.333000000001  # I137
$ws.Cells.Item(137, 10).Value = 3143.2222  # J137
$ws.Cells.Item(137, 11).Value = 26451.999  # K137
$ws.Cells.Item(137, 12).Value = 9429.6666  # L137
$ws.Cells.Item(137, 13).Value = -23901.999  # M137
$ws.Cells.Item(137, 14).Value = -14529.6666  # N137
# Row 138
$ws.Cells.Item(138, 8).Value = 6075.84  # H138
$ws.Cells.Item(138, 10).Value = 6164.75  # J138
$ws.Cells.Item(138, 12).Value = 18494.25  # L138
$ws.Cells.Item(138, 14).Value = -28774.25  # N138
$ws = $wb.Worksheets.Item("ARM")
# Row 18
$ws.Cells.Item(18, 8).Value = 3013  # H18
$ws.Cells.Item(18, 10).Value = 3013  # J18
$ws.Cells.Item(18, 12).Value = 3013  # L18
$ws.Cells.Item(18, 14).Value = -3657  # N18
# Row 61
$ws.Cells.Item(61, 8).Value = 7234.8823  # H61
$ws.Cells.Item(61, 9).Value = 5700.5  # I61
$ws.Cells.Item(61, 10).Value = 8598.777  # J61
$ws.Cells.Item(61, 11).Value = 5700.5  # K61
$ws.Cells.Item(61, 12).Value = 8598.777  # L61
$ws.Cells.Item(61, 13).Value = -5488.5  # M61
$ws.Cells.Item(61, 14).Value = -9022.777  # N61
# Row 74
$ws.Cells.Item(74, 8).Value = 7582.25  # H74
$ws.Cells.Item(74, 9).Value = 9744.954  # I74
$ws.Cells.Item(74, 10).Value = 2824.3  # J74
$ws.Cells.Item(74, 11).Value = 9744.954  # K74
$ws.Cells.Item(74, 12).Value = 2824.3  # L74
$ws.Cells.Item(74, 13).Value = -8870.954  # M74
$ws.Cells.Item(74, 14).Value = -4572.3  # N74
# Row 77
$ws.Cells.Item(77, 8).Value = 7582.25  # H77
$ws.Cells.Item(77, 9).Value = 9744.954  # I77
$ws.Cells.Item(77, 10).Value = 2824.3  # J77
$ws.Cells.Item(77, 11).Value = 48724.77  # K77
$ws.Cells.Item(77, 12).Value = 14121.5  # L77
$ws.Cells.Item(77, 13).Value = -44356.77  # M77
$ws.Cells.Item(77, 14).Value = -22857.5  # N77
# Row 122
$ws.Cells.Item(122, 8).Value = 9796.947  # H122
$ws.Cells.Item(122, 9).Value = 4445.4  # I122
$ws.Cells.Item(122, 10).Value = 15743.111  # J122
$ws.Cells.Item(122, 11).Value = 13336.2  # K122
$ws.Cells.Item(122, 12).Value = 47229.333  # L122
$ws.Cells.Item(122, 13).Value = -10886.2  # M122
$ws.Cells.Item(122, 14).Value = -52129.333  # N122
# Row 136
$ws.Cells.Item(136, 8).Value = 7234.8823  # H136
$ws.Cells.Item(136, 9).Value = 5700.5  # I136
$ws.Cells.Item(136, 10).Value = 8598.777  # J136
$ws.Cells.Item(136, 11).Value = 17101.5  # K136
$ws.Cells.Item(136, 12).Value = 25796.331  # L136
$ws.Cells.Item(136, 13).Value = -14551.5  # M136
$ws.Cells.Item(136, 14).Value = -30896.331  # N136
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Cells.Item(20, 8).Value = 2648.6287  # H20
$ws.Cells.Item(20, 9).Value = 1558.591  # I20
$ws.Cells.Item(20, 10).Value = 4493.3076  # J20
$ws.Cells.Item(20, 11).Value = 1558.591  # K20
$ws.Cells.Item(20, 12).Value = 4493.3076  # L20
$ws.Cells.Item(20, 13).Value = -1311.591  # M20
$ws.Cells.Item(20, 14).Value = -4987.3076  # N20
# Row 105
$ws.Cells.Item(105, 8).Value = 2987.5  # H105
$ws.Cells.Item(105, 9).Value = 2987.5  # I105
$ws.Cells.Item(105, 10).Value = 0  # J105
$ws.Cells.Item(105, 11).Value = 2987.5  # K105
$ws.Cells.Item(105, 12).Value = 0  # L105
$ws.Cells.Item(105, 13).ClearContents()  # M105
$ws.Cells.Item(105, 14).Value = -1240.5  # N105
# Row 124
$ws.Cells.Item(124, 8).Value = 0  # H124
$ws.Cells.Item(124, 10).Value = 0  # J124
$ws.Cells.Item(124, 12).ClearContents()  # L124
$ws.Cells.Item(124, 14).Value = 0  # N124
# Row 135
$ws.Cells.Item(135, 8).Value = 0  # H135
$ws.Cells.Item(135, 10).Value = 0  # J135
$ws.Cells.Item(135, 12).ClearContents()  # L135
$ws.Cells.Item(135, 14).Value = 0  # N135
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 44909276  # H31
$ws.Cells.Item(31, 9).Value = 12824320  # I31
$ws.Cells.Item(31, 10).Value = 76994230  # J31
$ws.Cells.Item(31, 11).Value = 12824320  # K31
$ws.Cells.Item(31, 12).Value = 76994230  # L31
$ws.Cells.Item(31, 13).Value = -12824025  # M31
$ws.Cells.Item(31, 14).Value = -76994820  # N31
# Row 34
$ws.Cells.Item(34, 8).Value = 44909276  # H34
$ws.Cells.Item(34, 9).Value = 12824320  # I34
$ws.Cells.Item(34, 10).Value = 76994230  # J34
$ws.Cells.Item(34, 11).Value = 12824320  # K34
$ws.Cells.Item(34, 12).Value = 76994230  # L34
$ws.Cells.Item(34, 13).Value = -12824118  # M34
$ws.Cells.Item(34, 14).Value = -76994634  # N34
# Row 62
$ws.Cells.Item(62, 8).Value = 99000  # H62
$ws.Cells.Item(62, 9).Value = 99000  # I62
$ws.Cells.Item(62, 11).Value = 99000  # K62
$ws.Cells.Item(62, 13).Value = -98376  # M62
# Row 65
$ws.Cells.Item(65, 8).Value = 99000  # H65
$ws.Cells.Item(65, 9).Value = 99000  # I65
$ws.Cells.Item(65, 11).Value = 495000  # K65
$ws.Cells.Item(65, 13).Value = -491880  # M65
# Row 141
$ws.Cells.Item(141, 8).Value = 168170.98  # H141
$ws.Cells.Item(141, 10).Value = 168170.98  # J141
$ws.Cells.Item(141, 12).Value = 168170.98  # L141
$ws.Cells.Item(141, 14).Value = -178530.98  # N141
$ws = $wb.Worksheets.Item("CUL")
# Row 124
$ws.Cells.Item(124, 8).Value = 3847.5  # H124
$ws.Cells.Item(124, 10).Value = 5479.5  # J124
$ws.Cells.Item(124, 12).Value = 16438.5  # L124
$ws.Cells.Item(124, 14).Value = -26258.5  # N124
# Row 131
$ws.Cells.Item(131, 8).Value = 13890301  # H131
$ws.Cells.Item(131, 9).Value = 166667500  # I131
$ws.Cells.Item(131, 10).Value = 1465.0454  # J131
$ws.Cells.Item(131, 11).Value = 500002500  # K131
$ws.Cells.Item(131, 12).Value = 4395.1362  # L131
$ws.Cells.Item(131, 13).Value = -499997460  # M131
$ws.Cells.Item(131, 14).Value = -14475.1362  # N131
# Row 133
$ws.Cells.Item(133, 8).Value = 5796.6  # H133
$ws.Cells.Item(133, 9).Value = 2994.3333  # I133
$ws.Cells.Item(133, 10).Value = 10000  # J133
$ws.Cells.Item(133, 11).Value = 8982.999899999999  # K133
$ws.Cells.Item(133, 12).Value = 30000  # L133
$ws.Cells.Item(133, 13).Value = -3922.999899999999  # M133
$ws.Cells.Item(133, 14).Value = -40120  # N133
# Row 134
$ws.Cells.Item(134, 8).Value = 5407.091  # H134
$ws.Cells.Item(134, 9).Value = 6568.2856  # I134
$ws.Cells.Item(134, 11).Value = 19704.8568  # K134
$ws.Cells.Item(134, 13).Value = -14634.8568  # M134
# Row 136
$ws.Cells.Item(136, 8).Value = 1579.122  # H136
$ws.Cells.Item(136, 9).Value = 1526.25  # I136
$ws.Cells.Item(136, 11).Value = 4578.75  # K136
$ws.Cells.Item(136, 13).Value = 521.25  # M136
# Row 137
$ws.Cells.Item(137, 8).Value = 51285892  # H137
$ws.Cells.Item(137, 9).Value = 4735.75  # I137
$ws.Cells.Item(137, 10).Value = 133335736  # J137
$ws.Cells.Item(137, 11).Value = 14207.25  # K137
$ws.Cells.Item(137, 12).Value = 400007208  # L137
$ws.Cells.Item(137, 13).Value = -9107.25  # M137
$ws.Cells.Item(137, 14).Value = -400017408  # N137
# Row 138
$ws.Cells.Item(138, 8).Value = 5493.3335  # H138
$ws.Cells.Item(138, 9).Value = 5493.3335  # I138
$ws.Cells.Item(138, 11).Value = 16480.0005  # K138
$ws.Cells.Item(138, 13).Value = -11340.0005  # M138
# Row 139
$ws.Cells.Item(139, 8).Value = 2087.5  # H139
$ws.Cells.Item(139, 9).Value = 2125.1  # I139
$ws.Cells.Item(139, 10).Value = 1899.5  # J139
$ws.Cells.Item(139, 11).Value = 6375.299999999999  # K139
$ws.Cells.Item(139, 12).Value = 5698.5  # L139
$ws.Cells.Item(139, 13).Value = -1235.299999999999  # M139
$ws.Cells.Item(139, 14).Value = -15978.5  # N139
# Row 140
$ws.Cells.Item(140, 8).Value = 6251574.5  # H140
$ws.Cells.Item(140, 9).Value = 17858170  # I140
$ws.Cells.Item(140, 11).Value = 53574510  # K140
$ws.Cells.Item(140, 13).Value = -53569330  # M140
# Row 141
$ws.Cells.Item(141, 8).Value = 3817.1  # H141
$ws.Cells.Item(141, 9).Value = 3817.1  # I141
$ws.Cells.Item(141, 11).Value = 11451.3  # K141
$ws.Cells.Item(141, 13).Value = -6271.299999999999  # M141
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Cells.Item(70, 8).Value = 62505812  # H70
$ws.Cells.Item(70, 9).Value = 4665.5557  # I70
$ws.Cells.Item(70, 10).Value = 142864430  # J70
$ws.Cells.Item(70, 11).Value = 4665.5557  # K70
$ws.Cells.Item(70, 12).Value = 142864430  # L70
$ws.Cells.Item(70, 13).Value = -4395.5557  # M70
$ws.Cells.Item(70, 14).Value = -142864970  # N70
# Row 73
$ws.Cells.Item(73, 8).Value = 62505812  # H73
$ws.Cells.Item(73, 9).Value = 4665.5557  # I73
$ws.Cells.Item(73, 10).Value = 142864430  # J73
$ws.Cells.Item(73, 11).Value = 4665.5557  # K73
$ws.Cells.Item(73, 12).Value = 142864430  # L73
$ws.Cells.Item(73, 13).Value = -3729.5557  # M73
$ws.Cells.Item(73, 14).Value = -142866302  # N73
# Row 132
$ws.Cells.Item(132, 8).Value = 45140.56  # H132
$ws.Cells.Item(132, 9).Value = 89426.336  # I132
$ws.Cells.Item(132, 10).Value = 4261.385  # J132
$ws.Cells.Item(132, 11).Value = 268279.008  # K132
$ws.Cells.Item(132, 12).Value = 12784.155  # L132
$ws.Cells.Item(132, 13).Value = -265749.008  # M132
$ws.Cells.Item(132, 14).Value = -17844.155  # N132
$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Cells.Item(122, 8).Value = 9222.666999999999  # H122
$ws.Cells.Item(122, 9).Value = 9449.833000000001  # I122
$ws.Cells.Item(122, 10).Value = 8768.333000000001  # J122
$ws.Cells.Item(122, 11).Value = 28349.499  # K122
$ws.Cells.Item(122, 12).Value = 26304.999  # L122
$ws.Cells.Item(122, 13).Value = -25899.499  # M122
$ws.Cells.Item(122, 14).Value = -31204.999  # N122
$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Cells.Item(54, 8).Value = 26795.4  # H54
$ws.Cells.Item(54, 10).Value = 19989.5  # J54
$ws.Cells.Item(54, 12).Value = 19989.5  # L54
$ws.Cells.Item(54, 14).Value = -21029.5  # N54
# Row 62
$ws.Cells.Item(62, 8).Value = 25012188  # H62
$ws.Cells.Item(62, 10).Value = 27789876  # J62
$ws.Cells.Item(62, 12).Value = 27789876  # L62
$ws.Cells.Item(62, 14).Value = -27791124  # N62
# Row 65
$ws.Cells.Item(65, 8).Value = 25012188  # H65
$ws.Cells.Item(65, 10).Value = 27789876  # J65
$ws.Cells.Item(65, 12).Value = 138949380  # L65
$ws.Cells.Item(65, 14).Value = -138955620  # N65
# Row 126
$ws.Cells.Item(126, 8).Value = 2277.0833  # H126
$ws.Cells.Item(126, 9).Value = 1675.7646  # I126
$ws.Cells.Item(126, 10).Value = 12499.5  # J126
$ws.Cells.Item(126, 11).Value = 5027.293799999999  # K126
$ws.Cells.Item(126, 12).Value = 37498.5  # L126
$ws.Cells.Item(126, 13).Value = -2557.293799999999  # M126
$ws.Cells.Item(126, 14).Value = -42438.5  # N126
# Row 129
$ws.Cells.Item(129, 8).Value = 20000  # H129
$ws.Cells.Item(129, 10).Value = 0  # J129
$ws.Cells.Item(129, 12).Value = 0  # L129
$ws.Cells.Item(129, 14).ClearContents()  # N129
# Row 131
$ws.Cells.Item(131, 8).Value = 58996.668  # H131
$ws.Cells.Item(131, 10).Value = 58996.668  # J131
$ws.Cells.Item(131, 12).Value = 58996.668  # L131
$ws.Cells.Item(131, 14).Value = -69076.66800000001  # N131
